$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting of the last two existing rows (88:89) onto the
# two new rows (90:91) so the new cells inherit the same cell styles.
$ws.Rows("88:89").Copy()
$ws.Rows("90:91").PasteSpecial()

# Row 90 - 四方坪站 (2025-10-15)
$ws.Range("A90").Value = 45945
$ws.Range("B90").Value = "四方坪站"
$ws.Range("C90").Formula = "=15812/126"
$ws.Range("D90").Formula = "=C90/(24*60)"
$ws.Range("E90").Formula = "=8539.95/126"
$ws.Range("F90").Formula = "=3006.5/126"
$ws.Range("G90").Formula = "=8539.95/(15812/60)"
$ws.Range("H90").Formula = "=375/126"

# Row 91 - 高岭站 (2025-10-15)
$ws.Range("A91").Value = 45945
$ws.Range("B91").Value = "高岭站"
$ws.Range("C91").Formula = "=6187/36"
$ws.Range("D91").Formula = "=C91/(24*60)"
$ws.Range("E91").Formula = "=4360.35/36"
$ws.Range("F91").Formula = "=1112.59/36"
$ws.Range("G91").Formula = "=4360.35/(6187/60)"
$ws.Range("H91").Formula = "=152/36"

# Update the view: selection moves to I90 and the frozen top-left scroll
# position is cleared (matches the post-edit sheetView).
$ws.Range("I90").Select()
